$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2466.9
$ws.Range("I32").Value = 2117.5
$ws.Range("J32").Value = 2699.8333
$ws.Range("K32").Value = 2117.5
$ws.Range("L32").Value = 2699.8333
$ws.Range("M32").Value = -1791.5
$ws.Range("N32").Value = -3351.8333
$ws.Range("H98").Value = 2433.182
$ws.Range("I98").Value = 1610.2667
$ws.Range("J98").Value = 4196.5713
$ws.Range("K98").Value = 1610.2667
$ws.Range("L98").Value = 4196.5713
$ws.Range("M98").Value = -112.2666999999999
$ws.Range("N98").Value = -7192.5713
$ws.Range("H113").Value = 4017.95
$ws.Range("I113").Value = 3843.9
$ws.Range("J113").Value = 4192
$ws.Range("K113").Value = 3843.9
$ws.Range("L113").Value = 4192
$ws.Range("M113").Value = -589.9000000000001
$ws.Range("N113").Value = -10700
$ws.Range("H122").Value = 2433.182
$ws.Range("I122").Value = 1610.2667
$ws.Range("J122").Value = 4196.5713
$ws.Range("K122").Value = 4830.800099999999
$ws.Range("L122").Value = 12589.7139
$ws.Range("M122").Value = -2380.800099999999
$ws.Range("N122").Value = -17489.7139
$ws.Range("H132").Value = 3705174.5
$ws.Range("I132").Value = 4445467
$ws.Range("J132").Value = 3710.6667
$ws.Range("K132").Value = 13336401
$ws.Range("L132").Value = 11132.0001
$ws.Range("M132").Value = -13333871
$ws.Range("N132").Value = -16192.0001
$ws.Range("H138").Value = 4011.2
$ws.Range("I138").Value = 2394.88
$ws.Range("J138").Value = 6031.6
$ws.Range("K138").Value = 7184.64
$ws.Range("L138").Value = 18094.8
$ws.Range("M138").Value = -2044.64
$ws.Range("N138").Value = -28374.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H21").Value = 32646.166
$ws.Range("I21").Value = 13960.75
$ws.Range("K21").Value = 13960.75
$ws.Range("M21").Value = -13586.75
$ws.Range("H32").Value = 3743.247
$ws.Range("I32").Value = 2414.4385
$ws.Range("K32").Value = 2414.4385
$ws.Range("M32").Value = -2127.4385
$ws.Range("H45").Value = 1907.76
$ws.Range("I45").Value = 984.2105
$ws.Range("K45").Value = 984.2105
$ws.Range("M45").Value = -607.2105
$ws.Range("H63").Value = 4900
$ws.Range("I63").Value = 2466.6667
$ws.Range("K63").Value = 2466.6667
$ws.Range("M63").Value = -1780.6667
$ws.Range("H66").Value = 4900
$ws.Range("I66").Value = 2466.6667
$ws.Range("K66").Value = 12333.3335
$ws.Range("M66").Value = -8901.333500000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H17").Value = 2504.5
$ws.Range("J17").Value = 2504.5
$ws.Range("L17").Value = 2504.5
$ws.Range("N17").Value = -2848.5
$ws.Range("H37").Value = 3596.5
$ws.Range("I37").Value = 2000
$ws.Range("J37").Value = 4128.6665
$ws.Range("K37").Value = 2000
$ws.Range("L37").Value = 4128.6665
$ws.Range("M37").Value = -1863
$ws.Range("N37").Value = -4402.6665

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 12962.692
$ws.Range("I10").Value = 437.375
$ws.Range("K10").Value = 437.375
$ws.Range("M10").Value = -298.375
$ws.Range("H12").Value = 10142815
$ws.Range("I12").Value = 13500418
$ws.Range("K12").Value = 13500418
$ws.Range("M12").Value = -13500248
$ws.Range("H31").Value = 2442141.5
$ws.Range("I31").Value = 4547389
$ws.Range("J31").Value = 4486.684
$ws.Range("K31").Value = 4547389
$ws.Range("L31").Value = 4486.684
$ws.Range("M31").Value = -4547094
$ws.Range("N31").Value = -5076.684
$ws.Range("H34").Value = 2442141.5
$ws.Range("I34").Value = 4547389
$ws.Range("J34").Value = 4486.684
$ws.Range("K34").Value = 4547389
$ws.Range("L34").Value = 4486.684
$ws.Range("M34").Value = -4547187
$ws.Range("N34").Value = -4890.684
$ws.Range("H86").Value = 3378
$ws.Range("I86").Value = 2435.7058
$ws.Range("J86").Value = 4445.933
$ws.Range("K86").Value = 2435.7058
$ws.Range("L86").Value = 4445.933
$ws.Range("M86").Value = -1312.7058
$ws.Range("N86").Value = -6691.933
$ws.Range("H89").Value = 3378
$ws.Range("I89").Value = 2435.7058
$ws.Range("J89").Value = 4445.933
$ws.Range("K89").Value = 12178.529
$ws.Range("L89").Value = 22229.665
$ws.Range("M89").Value = -6562.529
$ws.Range("N89").Value = -33461.665
$ws.Range("H132").Value = 2766.149
$ws.Range("I132").Value = 1749.2963
$ws.Range("K132").Value = 5247.8889
$ws.Range("M132").Value = -2717.8889
$ws.Range("H134").Value = 1659.9166
$ws.Range("I134").Value = 1246.6786
$ws.Range("K134").Value = 3740.0358
$ws.Range("M134").Value = -1205.0358

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 1599.091
$ws.Range("J47").Value = 3280.8
$ws.Range("L47").Value = 9842.400000000001
$ws.Range("N47").Value = -10704.4
$ws.Range("H68").Value = 1904.7826
$ws.Range("I68").Value = 727.3333
$ws.Range("J68").Value = 2532.7556
$ws.Range("K68").Value = 2181.9999
$ws.Range("L68").Value = 7598.266799999999
$ws.Range("M68").Value = -1370.9999
$ws.Range("N68").Value = -9220.266799999999
$ws.Range("H69").Value = 56409.55
$ws.Range("I69").Value = 581.1
$ws.Range("J69").Value = 112238
$ws.Range("K69").Value = 1743.3
$ws.Range("L69").Value = 336714
$ws.Range("M69").Value = -932.3000000000002
$ws.Range("N69").Value = -338336
$ws.Range("H71").Value = 1904.7826
$ws.Range("I71").Value = 727.3333
$ws.Range("J71").Value = 2532.7556
$ws.Range("K71").Value = 6545.9997
$ws.Range("L71").Value = 22794.8004
$ws.Range("M71").Value = -2489.9997
$ws.Range("N71").Value = -30906.8004
$ws.Range("H72").Value = 56409.55
$ws.Range("I72").Value = 581.1
$ws.Range("J72").Value = 112238
$ws.Range("K72").Value = 5229.900000000001
$ws.Range("L72").Value = 1010142
$ws.Range("M72").Value = -1173.900000000001
$ws.Range("N72").Value = -1018254

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 46819.668
$ws.Range("H16").Value = 46819.668
$ws.Range("H31").Value = 1000
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 1000
$ws.Range("M31").Value = -708
$ws.Range("H37").Value = 1000
$ws.Range("I37").Value = 1000
$ws.Range("K37").Value = 1000
$ws.Range("M37").Value = -723
$ws.Range("H59").Value = 15000
$ws.Range("J59").Value = 15000
$ws.Range("L59").Value = 15000
$ws.Range("N59").Value = -16166
$ws.Range("H70").Value = 3860
$ws.Range("I70").Value = 3800
$ws.Range("J70").Value = 3980
$ws.Range("K70").Value = 3800
$ws.Range("L70").Value = 3980
$ws.Range("M70").Value = -3530
$ws.Range("N70").Value = -4520
$ws.Range("H73").Value = 3860
$ws.Range("I73").Value = 3800
$ws.Range("J73").Value = 3980
$ws.Range("K73").Value = 3800
$ws.Range("L73").Value = 3980
$ws.Range("M73").Value = -2864
$ws.Range("N73").Value = -5852
$ws.Range("H97").Value = 1430.909
$ws.Range("I97").Value = 1294.4615
$ws.Range("J97").Value = 1628
$ws.Range("K97").Value = 1294.4615
$ws.Range("L97").Value = 1628
$ws.Range("M97").Value = -798.4614999999999
$ws.Range("N97").Value = -2620
$ws.Range("H107").Value = 669.75
$ws.Range("I107").Value = 190.625
$ws.Range("J107").Value = 989.1667
$ws.Range("K107").Value = 190.625
$ws.Range("L107").Value = 989.1667
$ws.Range("M107").Value = 1729.375
$ws.Range("N107").Value = -4829.1667
$ws.Range("H125").Value = 24375
$ws.Range("J125").Value = 24375
$ws.Range("L125").Value = 24375
$ws.Range("N125").Value = -29295

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 2216.375
$ws.Range("I35").Value = 1390.1428
$ws.Range("J35").Value = 8000
$ws.Range("K35").Value = 1390.1428
$ws.Range("L35").Value = 8000
$ws.Range("M35").Value = -1054.1428
$ws.Range("N35").Value = -8672
$ws.Range("H40").Value = 2446.4614
$ws.Range("I40").Value = 1766.3334
$ws.Range("J40").Value = 2650.5
$ws.Range("K40").Value = 1766.3334
$ws.Range("L40").Value = 2650.5
$ws.Range("M40").Value = -1630.3334
$ws.Range("N40").Value = -2922.5
$ws.Range("H45").Value = 7126
$ws.Range("I45").Value = 3941
$ws.Range("K45").Value = 3941
$ws.Range("M45").Value = -3534
$ws.Range("H106").Value = 13450
$ws.Range("J106").Value = 13450
$ws.Range("L106").Value = 13450
$ws.Range("N106").Value = -15974
$ws.Range("H136").Value = 1924628.5
$ws.Range("I136").Value = 2501197
$ws.Range("K136").Value = 7503591
$ws.Range("M136").Value = -7501041

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 61670.668
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()
$ws.Range("H58").Value = 14934.643
$ws.Range("I58").Value = 14085
$ws.Range("K58").Value = 14085
$ws.Range("M58").Value = -13777
$ws.Range("H113").Value = 1644.2
$ws.Range("I113").Value = 524.25
$ws.Range("J113").Value = 2924.1428
$ws.Range("K113").Value = 1572.75
$ws.Range("L113").Value = 8772.428400000001
$ws.Range("M113").Value = 597.25
$ws.Range("N113").Value = -13112.4284
$ws.Range("H122").Value = 478741.8
$ws.Range("I122").Value = 668325.2
$ws.Range("J122").Value = 4783.3335
$ws.Range("K122").Value = 2004975.6
$ws.Range("L122").Value = 14350.0005
$ws.Range("M122").Value = -2002525.6
$ws.Range("N122").Value = -19250.0005
$ws.Range("H123").Value = 65000
$ws.Range("J123").Value = 65000
$ws.Range("L123").Value = 65000
$ws.Range("N123").Value = -74800
$ws.Range("H132").Value = 182252.14
$ws.Range("I132").Value = 239410.05
$ws.Range("J132").Value = 10778.429
$ws.Range("K132").Value = 718230.1499999999
$ws.Range("L132").Value = 32335.287
$ws.Range("M132").Value = -715700.1499999999
$ws.Range("N132").Value = -37395.287
$ws.Range("H136").Value = 1044.5593
$ws.Range("I136").Value = 543.9434
$ws.Range("K136").Value = 1631.8302
$ws.Range("M136").Value = 918.1698000000001
